# Applies crypto price/volume updates scraped on Wed Oct  9 16:46:22 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. '578.56' or thousands-dotted
# '61.805.55'); format as Text first so Excel keeps them as strings instead
# of auto-converting to numbers, then restore the default style so no
# formatting footprint is left behind.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "61.805.55"
$ws.Range("E2").Value = "  -1.04%  "
Set-TextValue $ws.Range("D3") "2.445.48"
$ws.Range("E3").Value = "  -0.17%  "
Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws.Range("D5") "578.56"
$ws.Range("E5").Value = "  -0.37%  "
Set-TextValue $ws.Range("D6") "140.79"
$ws.Range("E6").Value = "  -1.99%  "
$ws.Range("E7").Value = "  +0.12%  "
Set-TextValue $ws.Range("D8") "0.532"
$ws.Range("E8").Value = "  +0.60%  "
Set-TextValue $ws.Range("D9") "2.436.18"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("E13").Value = "  -2.25%  "
Set-TextValue $ws.Range("D14") "25.90"
Set-TextValue $ws.Range("D15") "2.906.39"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("E16").Value = "  -1.27%  "
Set-TextValue $ws.Range("D17") "61.729.61"
Set-TextValue $ws.Range("D18") "2.445.69"
$ws.Range("E18").Value = "  +0.83%  "
Set-TextValue $ws.Range("D19") "10.61"
$ws.Range("E19").Value = "  -3.82%  "
$ws.Range("E20").Value = "  +1.55%  "
Set-TextValue $ws.Range("D21") "324.87"
$ws.Range("E21").Value = "  -2.39%  "
Set-TextValue $ws.Range("D22") "4.08"
$ws.Range("E22").Value = "  -1.35%  "
Set-TextValue $ws.Range("D23") "6.01"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("E25").Value = "  +0.01%  "
Set-TextValue $ws.Range("D26") "64.96"
$ws.Range("E26").Value = "  -1.48%  "
Set-TextValue $ws.Range("D27") "9.08"
$ws.Range("E27").Value = "  -0.76%  "
Set-TextValue $ws.Range("D28") "582.14"
$ws.Range("E28").Value = "  -8.92%  "
Set-TextValue $ws.Range("D29") "2.572.99"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  -3.77%  "
Set-TextValue $ws.Range("D32") "7.93"
$ws.Range("E32").Value = "  -1.94%  "
Set-TextValue $ws.Range("D33") "1.36"
$ws.Range("E33").Value = "  -5.97%  "
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("E35").Value = "  -5.65%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  -5.42%  "
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D39") "1.40"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D40") "150.97"
$ws.Range("E40").Value = "  +1.57%  "
Set-TextValue $ws.Range("D41") "18.26"
$ws.Range("E41").Value = "  -0.92%  "
Set-TextValue $ws.Range("D42") "5.14"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("E45").Value = "  -2.54%  "
Set-TextValue $ws.Range("D46") "2.37"
$ws.Range("E46").Value = "  -5.61%  "
$ws.Range("E47").Value = "  +22.65%  "
Set-TextValue $ws.Range("D48") "142.76"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("E49").Value = "  -2.68%  "
Set-TextValue $ws.Range("D50") "0.598"
$ws.Range("E50").Value = "  -0.14%  "
Set-TextValue $ws.Range("D51") "19.69"
$ws.Range("E51").Value = "  -0.46%  "
